# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: volume number + report week dates ----
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# ---- Crime Complaints table (rows 14-31) ----
# Row 14: Murder
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -33.333333333333
$ws.Range("I14").Value = 55
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -8.333333333333
$ws.Range("M14").Value = -51.327433628318
$ws.Range("N14").Value = -85.488126649076

# Row 15: Rape
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 33.333333333333
$ws.Range("F15").Value = 23
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 43.75
$ws.Range("I15").Value = 198
$ws.Range("J15").Value = 169
$ws.Range("K15").Value = 17.159763313609
$ws.Range("L15").Value = -0.502512562814
$ws.Range("M15").Value = 22.222222222222
$ws.Range("N15").Value = -57.782515991471

# Row 16: Robbery
$ws.Range("C16").Value = 37
$ws.Range("D16").Value = 44
$ws.Range("E16").Value = -15.909090909090
$ws.Range("F16").Value = 183
$ws.Range("G16").Value = 190
$ws.Range("H16").Value = -3.684210526315
$ws.Range("I16").Value = 1835
$ws.Range("J16").Value = 1870
$ws.Range("K16").Value = -1.871657754010
$ws.Range("L16").Value = -8.387418871692
$ws.Range("M16").Value = -31.784386617100
$ws.Range("N16").Value = -85.466497703152

# Row 17: Fel. Assault
$ws.Range("C17").Value = 83
$ws.Range("D17").Value = 90
$ws.Range("E17").Value = -7.777777777777
$ws.Range("F17").Value = 329
$ws.Range("G17").Value = 355
$ws.Range("H17").Value = -7.323943661971
$ws.Range("I17").Value = 3372
$ws.Range("J17").Value = 3314
$ws.Range("K17").Value = 1.750150875075
$ws.Range("L17").Value = 4.785581106277
$ws.Range("M17").Value = 29.344073647871
$ws.Range("N17").Value = -49.422528873556

# Row 18: Burglary
$ws.Range("C18").Value = 41
$ws.Range("D18").Value = 44
$ws.Range("E18").Value = -6.818181818181
$ws.Range("F18").Value = 147
$ws.Range("G18").Value = 159
$ws.Range("H18").Value = -7.547169811320
$ws.Range("I18").Value = 1466
$ws.Range("J18").Value = 1595
$ws.Range("K18").Value = -8.087774294670
$ws.Range("L18").Value = -20.195971692977
$ws.Range("M18").Value = -39.769926047658
$ws.Range("N18").Value = -84.246722544594

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 111
$ws.Range("D19").Value = 116
$ws.Range("E19").Value = -4.310344827586
$ws.Range("F19").Value = 421
$ws.Range("G19").Value = 460
$ws.Range("H19").Value = -8.478260869565
$ws.Range("I19").Value = 4079
$ws.Range("J19").Value = 4482
$ws.Range("K19").Value = -8.991521642124
$ws.Range("L19").Value = -10.292500549813
$ws.Range("M19").Value = 26.128633271490
$ws.Range("N19").Value = -22.804693414080

# Row 20: G.L.A.
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 38
$ws.Range("E20").Value = 21.052631578947
$ws.Range("F20").Value = 194
$ws.Range("G20").Value = 167
$ws.Range("H20").Value = 16.167664670658
$ws.Range("I20").Value = 1380
$ws.Range("J20").Value = 1407
$ws.Range("K20").Value = -1.918976545842
$ws.Range("L20").Value = -1.988636363636
$ws.Range("M20").Value = 26.489459211732
$ws.Range("N20").Value = -80.812013348164

# Row 21: TOTAL
$ws.Range("C21").Value = 323
$ws.Range("D21").Value = 337
$ws.Range("E21").Value = -4.154302670623
$ws.Range("F21").Value = 1301
$ws.Range("G21").Value = 1353
$ws.Range("H21").Value = -3.843311160384
$ws.Range("I21").Value = 12385
$ws.Range("J21").Value = 12892
$ws.Range("K21").Value = -3.932671424139
$ws.Range("L21").Value = -6.683242917420
$ws.Range("M21").Value = 0.437920687697
$ws.Range("N21").Value = -70.457743959163

# Row 22: Transit
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 5.263157894736
$ws.Range("I22").Value = 219
$ws.Range("J22").Value = 220
$ws.Range("K22").Value = -0.454545454545
$ws.Range("L22").Value = -18.587360594795
$ws.Range("M22").Value = -31.775700934579
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 36
$ws.Range("E23").Value = -22.222222222222
$ws.Range("F23").Value = 93
$ws.Range("G23").Value = 116
$ws.Range("H23").Value = -19.827586206896
$ws.Range("I23").Value = 1121
$ws.Range("J23").Value = 1222
$ws.Range("K23").Value = -8.265139116202
$ws.Range("L23").Value = -4.676870748299
$ws.Range("M23").Value = 23.594266813671
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 231
$ws.Range("D24").Value = 206
$ws.Range("E24").Value = 12.135922330097
$ws.Range("F24").Value = 1007
$ws.Range("G24").Value = 967
$ws.Range("H24").Value = 4.136504653567
$ws.Range("I24").Value = 9452
$ws.Range("J24").Value = 9563
$ws.Range("K24").Value = -1.160723622294
$ws.Range("L24").Value = -8.054474708171
$ws.Range("M24").Value = 19.027830248079
$ws.Range("N24").Value = "***.*"

# Row 25: Retail Theft
$ws.Range("C25").Value = 104
$ws.Range("D25").Value = 64
$ws.Range("E25").Value = 62.5
$ws.Range("F25").Value = 426
$ws.Range("G25").Value = 337
$ws.Range("H25").Value = 26.409495548961
$ws.Range("I25").Value = 4247
$ws.Range("J25").Value = 3801
$ws.Range("K25").Value = 11.733754275190
$ws.Range("L25").Value = -1.209583624098
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"

# Row 26: Misd. Assault
$ws.Range("C26").Value = 145
$ws.Range("D26").Value = 124
$ws.Range("E26").Value = 16.935483870967
$ws.Range("F26").Value = 540
$ws.Range("G26").Value = 466
$ws.Range("H26").Value = 15.879828326180
$ws.Range("I26").Value = 4912
$ws.Range("J26").Value = 4714
$ws.Range("K26").Value = 4.200254560882
$ws.Range("L26").Value = 7.719298245614
$ws.Range("M26").Value = -20.311486048020
$ws.Range("N26").Value = "***.*"

# Row 27: UCR Rape*
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 60
$ws.Range("F27").Value = 29
$ws.Range("G27").Value = 26
$ws.Range("H27").Value = 11.538461538461
$ws.Range("I27").Value = 275
$ws.Range("J27").Value = 262
$ws.Range("K27").Value = 4.961832061068
$ws.Range("L27").Value = -7.094594594594
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Other Sex Crimes
$ws.Range("C28").Value = 21
$ws.Range("D28").Value = 14
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 67
$ws.Range("G28").Value = 62
$ws.Range("H28").Value = 8.064516129032
$ws.Range("I28").Value = 513
$ws.Range("J28").Value = 497
$ws.Range("K28").Value = 3.219315895372
$ws.Range("L28").Value = 8.917197452229
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"

# Row 29: Shooting Vic.
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = -10
$ws.Range("F29").Value = 23
$ws.Range("G29").Value = 25
$ws.Range("H29").Value = -8
$ws.Range("I29").Value = 194
$ws.Range("J29").Value = 188
$ws.Range("K29").Value = 3.191489361702
$ws.Range("L29").Value = -29.454545454545
$ws.Range("M29").Value = -54.352941176470
$ws.Range("N29").Value = -87.160820648577

# Row 30: Shooting Inc.
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = 16.666666666666
$ws.Range("F30").Value = 15
$ws.Range("G30").Value = 17
$ws.Range("H30").Value = -11.764705882352
$ws.Range("I30").Value = 158
$ws.Range("J30").Value = 158
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -30.701754385964
$ws.Range("M30").Value = -54.069767441860
$ws.Range("N30").Value = -88.382352941176

# Row 31: Hate Crimes
$ws.Range("C31").Value = "0"
$ws.Range("D31").Value = "0"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 63
$ws.Range("J31").Value = 53
$ws.Range("K31").Value = 18.867924528301
$ws.Range("L31").Value = -5.970149253731
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"

# ---- Traffic Statistics (row 33: Traffic Fatalities) ----
$ws.Range("G33").Value = "0"
$ws.Range("H33").Value = "***.*"
$ws.Range("L33").Value = -13.043478260869

# ---- Column E width shrinks slightly as the widest value in the column got shorter ----
# (now matches the width already used by columns C/D/F/G/H)
$ws.Columns("E").ColumnWidth = $ws.Columns("C").ColumnWidth
